$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '26.334.51'
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).Value = '  -4.03%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.766.09'
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = '  -2.98%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9935'
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = '  -1.00%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '314.64'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +0.81%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.9960'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -0.63%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.4271'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = '  +1.00%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3635'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = '  +1.32%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.07172'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +0.06%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.8483'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -0.38%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '20.32'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -0.22%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '1.763.31'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = '  -3.20%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '5.238'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -2.12%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Chainlink'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '6.420'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = '  +0.12%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.06836'
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = '  -1.46%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.9993'
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = '  -0.64%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '78.78'
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).Value = '  -3.25%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.000008667'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -2.28%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.9973'
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = '  -0.51%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '15.00'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -1.87%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '26.339.45'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -3.98%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '5.085'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -0.54%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '11.09'
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +1.46%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '1.983.87'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  -4.04%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '151.37'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -1.74%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '1.856'
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -6.30%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '18.02'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -1.48%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '5.087'
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = '  -0.43%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '113.62'
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = '  -0.31%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.788'
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +1.72%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.08923'
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = '  +0.22%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.7281'
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = '  -2.00%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.120'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +0.66%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.313'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -4.23%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.9921'
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -1.08%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.711'
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -8.24%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.083'
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = '  +0.36%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.05147'
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -1.41%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.01887'
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = '  -1.36%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.4915'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -2.45%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.1607'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = '  -2.49%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '2.581'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -6.87%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '6.300'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +0.16%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '8.003'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -3.26%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '10.24'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -1.71%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Quant'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '104.69'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = '  -1.10%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.9956'
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -0.67%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.4481'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -3.46%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.06172'
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -4.14%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.604'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -0.27%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.740'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +3.71%  '
